# MIAPPE_Observation_unit_and_sample.xlsx
# "added term accession numbers to ambiguous tags and harmonized similar tags"
#
# The "SwateTemplateMetadata" sheet keeps a #TAGS list in row 12 (Tags),
# row 13 (Tags Term Accession Number) and row 14 (Tags Term Source REF).
# The tag "Study" is harmonized to lower-case "study" (matching the style
# of the neighbouring "growth protocol" tag) and gets a proper term
# accession number + term source REF (NCIT) added, where before it had
# none.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Harmonize the "Study" tag to lower-case "study"
$ws.Range("D12").Value = "study"

# Add the missing term accession number (NCIT) for the "study" tag and
# turn it into a hyperlink, like the other accession-number cells.
$ws.Hyperlinks.Add($ws.Range("D13"), "http://purl.obolibrary.org/obo/NCIT_C63536")
$ws.Range("D13").Style = "Hyperlink"

# Add the matching term source REF for the "study" tag.
$ws.Range("D14").Value = "NCIT"

# The SwateTemplateMetadata sheet becomes the active/selected sheet and
# the newly-added accession number cell becomes the active selection.
$ws.Activate()
$ws.Range("D13").Select()
